# aggiornamento fino a 20/09/2021
# Append new daily rows (375-385) to the data table on Sheet1, reusing the
# formatting (date style) of the last existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 374
$firstNewRow = 375
$lastNewRow = 385

# Copy the formatting (number format / style / borders) of the last data
# row down into the newly appended rows before writing values into them.
$ws.Range("A$lastRow`:D$lastRow").Copy()
$ws.Range("A$firstNewRow`:D$lastNewRow").PasteSpecial(-4122)

$data = @(
    @(375, 44449, 4, 12, 106.2981663566303),
    @(376, 44450, 1, 13, 115.1563468863496),
    @(377, 44451, 3, 16, 141.7308884755071),
    @(378, 44452, 4, 20, 177.1636105943839),
    @(379, 44453, 2, 22, 194.8799716538223),
    @(380, 44454, 0, 19, 168.3054300646647),
    @(381, 44455, 3, 17, 150.5890690052263),
    @(382, 44456, 1, 14, 124.0145274160687),
    @(383, 44457, 6, 19, 168.3054300646647),
    @(384, 44458, 4, 20, 177.1636105943839),
    @(385, 44459, 0, 16, 141.7308884755071)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
